# chore(exportDoc): add role column to exported excel sheet
#
# Inserts a new "Role" column between "MiddleName" and "Branch" on the
# "Submitted Claims" sheet, shifting the existing Branch..Month of Claim
# columns one place to the right, and fills in the Role values for each
# staff row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F (pushes Branch..Month of Claim -> G..N)
$ws.Range("F1").EntireColumn.Insert()

# New header
$ws.Range("F1").Value = "Role"

# New per-row values
$ws.Range("F2").Value = "Service Executive (Financial)"
$ws.Range("F3").Value = "Service Executive (Financial)"
$ws.Range("F4").Value = "Service Executive (Non-Financial)"
